$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# DATE_TYPE_CODE (must remain text "001", not numeric 1)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# REPORT_DATE
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# PARENT_NETPROFIT
$ws.Range("O2").Value = 179826999.88

# TOTAL_OPERATE_INCOME
$ws.Range("P2").Value = 5311212819.9

# TOTAL_OPERATE_COST
$ws.Range("Q2").Value = 5126617692.36

# TOE_RATIO (was empty)
$ws.Range("R2").Value = 71.7571004062

# OPERATE_COST
$ws.Range("S2").Value = 4657329904.51

# OPERATE_EXPENSE
$ws.Range("T2").Value = 4657329904.51

# OPERATE_EXPENSE_RATIO (was empty)
$ws.Range("U2").Value = 71.3810659336

# SALE_EXPENSE
$ws.Range("V2").Value = 74306253.5

# MANAGE_EXPENSE
$ws.Range("W2").Value = 113713644.39

# FINANCE_EXPENSE
$ws.Range("X2").Value = 90270632.97

# OPERATE_PROFIT
$ws.Range("Y2").Value = 212309404.28

# TOTAL_PROFIT
$ws.Range("Z2").Value = 211928544.15

# INCOME_TAX
$ws.Range("AA2").Value = 32101544.27

# OPERATE_TAX_ADD
$ws.Range("AG2").Value = 15527044.96

# TOI_RATIO (was empty)
$ws.Range("AP2").Value = 73.10424485990001

# OPERATE_PROFIT_RATIO (was empty)
$ws.Range("AQ2").Value = 134.903766388341

# PARENT_NETPROFIT_RATIO (was empty)
$ws.Range("AR2").Value = 184.834165558612

# DEDUCT_PARENT_NETPROFIT
$ws.Range("AS2").Value = 121266643.64

# DPN_RATIO (was empty)
$ws.Range("AT2").Value = 176.740779630187
